$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price snapshot was scraped. Insert a fresh column at Y (pushing the
# existing "nom" column to Z and "url_produit" to AA), stamp the new column's
# header with the scrape timestamp, and seed every product row's new price
# with the same value as the previous (rightmost, column X) snapshot - the
# site wasn't re-scraped for this run, it's just a repeated/rolled-forward
# reading.

$newColIndex = 25   # column Y
$priceColIndex = 24 # column X (last existing snapshot column)

$ws.Columns($newColIndex).Insert()

$ws.Cells.Item(1, $newColIndex).Value = "2026-01-28 18:23:01"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $priceValue = $ws.Cells.Item($r, $priceColIndex).Value2
    $ws.Cells.Item($r, $newColIndex).Value = $priceValue
}
